$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric values in columns B..E (rows 2..13) to the nearest
# integer, matching the "write as integer data" change described in the
# commit message.
for ($r = 2; $r -le 13; $r++) {
    for ($c = 2; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = $excel.WorksheetFunction.Round([double]$cell.Value2, 0)
    }
}
